$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4391
$wsExhibit.Range("F12").Value = 1646
$wsExhibit.Range("F14").Value = 3506

# Sheet "全部类型" (fourth sheet)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4391
$wsAll.Range("F16").Value = 1646
$wsAll.Range("F18").Value = 3506

$wb.Save()
